# "changes for parallel execution"
# Adds 7 new "post order background" rows at the top of the data on the
# "Post" sheet (pushing the existing scenarios down), and updates the
# remembered cell-selection on the "Post" and "Delete" sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Post")
$ws4 = $wb.Worksheets.Item("Delete")

# Insert 7 blank rows above the current row 2 -- this pushes the existing
# scenario rows (2-6) down to rows 9-13 and keeps their values/styles intact.
$ws1.Rows("2:8").Insert() | Out-Null

$names = @("test name1", "test name2", "test name3", "test name4", "test name5", "test name6", "test name7")
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 2 + $i
    $ws1.Range("A$r").Value = "POST"
    $ws1.Range("B$r").Value = "post order background"
    $ws1.Range("C$r").Value = 1
    $ws1.Range("D$r").Value = $names[$i]
}

# Column E holds the literal text "true" (not the Boolean value). Typing the
# string directly auto-converts to a Boolean, so instead copy it down from
# the existing "true" cell (now at E9, post-insert) which preserves both the
# literal text type and the cell style.
for ($r = 2; $r -le 8; $r++) {
    $ws1.Range("E9").Copy($ws1.Range("E$r")) | Out-Null
}
$excel.CutCopyMode = 0

# Update the remembered selections.
$ws1.Range("C5").Select() | Out-Null
$ws4.Range("B3").Select() | Out-Null

# Leave "Post" as the active sheet/tab, as it was before the edit.
$ws1.Activate() | Out-Null
